$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 14: new weekly price entry (was dated 44475, now 44488) ---
$ws.Range("D14").Value = 44488
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 26000
$ws.Range("P14").Value = 25600
$ws.Range("S14").Value = 2560

# --- Update existing row 15: shift its date forward (was 44461, now 44475) ---
$ws.Range("D15").Value = 44475

# --- Insert a new row 16 containing the data that used to live in row 15 ---
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"

$ws.Range("D16").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("D16").Value = 44461

$ws.Range("E16").Value = 8
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107002
$ws.Range("J16").Value = "Chirimoya"
$ws.Range("K16").Value = "Cultivar IV Región"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 29000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 29500
$ws.Range("Q16").Value = "$/bandeja 10 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 2950
$ws.Range("T16").Value = 10
